$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (AD:AF) reusing the existing header style
# (bold, centered, thin-bordered) from the last header cell, then set text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row (2-50).
$ws.Range("AD2:AD50").Value = 63
$ws.Range("AE2:AE50").Value = 99
$ws.Range("AF2:AF50").Value = 0
